$wb = $excel.ActiveWorkbook

# --- 1. Update test-run timestamps embedded in the e-mail addresses on the
#        "UsuariosRegistro" sheet (120951 -> 130229) ---
$wsUsuarios = $wb.Worksheets.Item("UsuariosRegistro")
for ($r = 2; $r -le 6; $r++) {
    $cell = $wsUsuarios.Cells.Item($r, 3)
    $oldValue = $cell.Value()
    $cell.Value = $oldValue.Replace("20251110_120951", "20251110_130229")
}

# The "LoginData" sheet reuses the same Juan/María e-mail addresses (rows 2-3,
# column A) for its successful-login fixtures, so keep those in sync too.
$wsLogin = $wb.Worksheets.Item("LoginData")
for ($r = 2; $r -le 3; $r++) {
    $cell = $wsLogin.Cells.Item($r, 1)
    $oldValue = $cell.Value()
    $cell.Value = $oldValue.Replace("20251110_120951", "20251110_130229")
}

# --- 2. Replace the product catalog on the "ProductosBusqueda" sheet ---
$wsProductos = $wb.Worksheets.Item("ProductosBusqueda")

# Only touch the cells whose value actually changes so cells that stay the
# same (e.g. already-blank SubCategoria cells) keep their original state.
# A lone leading apostrophe is Excel's "treat as text" entry marker: it keeps
# the cell a text cell (matching the rest of this all-text sheet) instead of
# auto-converting to a number / going fully blank, while leaving the stored
# value itself empty, same as the other already-blank cells in this table.
$wsProductos.Cells.Item(2, 1).Value = "Laptops & Notebooks"
$wsProductos.Cells.Item(2, 2).Value = "'"
$wsProductos.Cells.Item(2, 3).Value = "MacBook"

$wsProductos.Cells.Item(3, 1).Value = "'"
$wsProductos.Cells.Item(3, 3).Value = "iPhone"
$wsProductos.Cells.Item(3, 4).Value = "'1"

$wsProductos.Cells.Item(4, 1).Value = "Cameras"
$wsProductos.Cells.Item(4, 2).Value = "'"
$wsProductos.Cells.Item(4, 3).Value = "Canon EOS 5D"

$wsProductos.Cells.Item(5, 1).Value = "Laptops & Notebooks"
$wsProductos.Cells.Item(5, 2).Value = "Macs"
$wsProductos.Cells.Item(5, 3).Value = "MacBook Air"
$wsProductos.Cells.Item(5, 4).Value = "'2"

$wsProductos.Cells.Item(6, 1).Value = "Tablets"
$wsProductos.Cells.Item(6, 3).Value = "Samsung Galaxy Tab 10.1"

# --- 3. Widen column C to fit the new, longer product names ---
# (Excel quantises ColumnWidth to whole pixels, so 20.3 is the input that
# lands closest to the target stored width of 21.24609375 characters.)
$wsProductos.Columns.Item(3).ColumnWidth = 20.3
